$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 12: "AutoBaseBuildingManager" ---
# (pushes the existing "Stack" row and everything below it down by one)
$ws.Rows(12).Insert()

# Copy formatting from row 3 (style used by "HashNode", a plain Critical/Classes/No/Zombie/N/A row)
$ws.Range("A3:G3").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)

$ws.Cells.Item(12,1).Value2 = "AutoBaseBuildingManager"
$ws.Cells.Item(12,2).Value2 = "Critical"
$ws.Cells.Item(12,3).Value2 = "CriticalScars\AutoBaseBuilder\Classes"
$ws.Cells.Item(12,4).Value2 = "No"
$ws.Cells.Item(12,5).Value2 = "Zombie"
$ws.Cells.Item(12,6).Value2 = "N/A"

# --- Insert new row 28: "SquadDataHistory" ---
# After the first insertion, "SquadData" (originally row 26) is now row 27.
# We insert directly below it.
$ws.Rows(28).Insert()

# Copy formatting from row 27 (now "SquadData", style used for Limited/Parent-type rows)
$ws.Range("A27:G27").Copy()
$ws.Range("A28:G28").PasteSpecial(-4122)

$ws.Cells.Item(28,1).Value2 = "SquadDataHistory"
$ws.Cells.Item(28,2).Value2 = "Critical"
$ws.Cells.Item(28,3).Value2 = "CriticalScars\SquadManager\Classes"
$ws.Cells.Item(28,4).Value2 = "Yes"
$ws.Cells.Item(28,5).Value2 = "Limited"
$ws.Cells.Item(28,6).Value2 = "SquadData, FixedLengthTable"

$excel.CutCopyMode = 0
